# Generate Report for Handback
# Updates priority (ht -> mt) and timestamps for the two files that were
# re-processed (453447bb... and b1289512...) across the Overview, zh-cn,
# and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
$wsOverview.Range("G2").Value = "2016-08-24 04:15:27"
$wsOverview.Range("G4").Value = "2016-08-24 04:15:27"

# --- zh-cn sheet ---
# Priority column (E): ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"

# Correspond Handoff Datetime column (H)
$wsZhCn.Range("H2").Value = "2016-08-24 04:15:22"
$wsZhCn.Range("H4").Value = "2016-08-24 04:15:22"

# Correspond Handback DateTime column (K)
$wsZhCn.Range("K2").Value = "2016-08-24 04:15:39"
$wsZhCn.Range("K4").Value = "2016-08-24 04:15:39"

# --- de-de sheet ---
# Priority column (E): ht -> mt
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# Correspond Handoff Datetime column (H)
$wsDeDe.Range("H2").Value = "2016-08-24 04:15:27"
$wsDeDe.Range("H4").Value = "2016-08-24 04:15:27"

# Correspond Handback DateTime column (K)
$wsDeDe.Range("K2").Value = "2016-08-24 04:15:46"
$wsDeDe.Range("K4").Value = "2016-08-24 04:15:46"
